$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10775351.92785637
$ws.Range("C2").Value = 2429237.833225261
$ws.Range("D2").Value = 28250397.67641928
$ws.Range("E2").Value = 1176877.84930499
$ws.Range("F2").Value = 9397656.956019122
$ws.Range("G2").Value = 1848931.834534268
$ws.Range("H2").Value = 2100844.41724581
$ws.Range("I2").Value = 10775351.92785637
$ws.Range("J2").Value = 44940247
$ws.Range("K2").Value = 112
$ws.Range("L2").Value = 30679635.50964454
$ws.Range("M2").Value = 10574534.80532411
$ws.Range("N2").Value = 3949776.251780078
$ws.Range("O2").Value = 41293.57878157049
$ws.Range("P2").Value = 206574.1005918898
$ws.Range("Q2").Value = 247867.6793734603
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 109448.5612180983
$ws.Range("T2").Value = 109448.5612180983
